# This script reproduces the "Updated cryptos list" data refresh.
# It rewrites the Coin / Link / Price / Volume(1h) cells that changed between
# runs (including two rows, 27/28 and 50/51, whose coins were re-ordered).
#
# Columns D (Price) and E (Volume) are stored as TEXT in the sheet (e.g. "1.00",
# "62.303.25", "  -2.25%  "). Assigning such look-alike-numeric strings straight
# to .Value would make Excel auto-convert them to real numbers, so for any new
# Price value that Excel would parse as a plain number we temporarily force the
# cell to Text format, assign the value, then restore the original ("Normal")
# cell style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.303.25"
$ws.Range("E2").Value = "  -2.25%  "

# Row 3
$ws.Range("D3").Value = "3.010.27"
$ws.Range("E3").Value = "  -2.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.79%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.30%  "

# Row 9
$ws.Range("D9").Value = "3.006.95"
$ws.Range("E9").Value = "  -2.15%  "

# Row 10
$ws.Range("E10").Value = "  -5.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.13%  "

# Row 13
$ws.Range("E13").Value = "  -4.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.64%  "

# Row 15
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").Value = "3.500.89"
$ws.Range("E16").Value = "  -2.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.06%  "

# Row 18
$ws.Range("D18").Value = "62.258.95"
$ws.Range("E18").Value = "  -2.22%  "

# Row 19
$ws.Range("D19").Value = "3.005.30"
$ws.Range("E19").Value = "  -2.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.62%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.60%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.07%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$ws.Range("E30").Value = "  -3.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.73%  "

# Row 34
$ws.Range("E34").Value = "  -3.41%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0797"
$ws.Range("E35").Value = "  -3.21%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.83%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "

# Row 41
$ws.Range("E41").Value = "  -13.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "390.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.83%  "

# Row 44
$ws.Range("E44").Value = "  -2.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.269"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.80%  "

# Row 46
$ws.Range("D46").Value = "2.725.37"
$ws.Range("E46").Value = "  -4.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.19%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.31"
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "  +0.08%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.20%  "
